$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Janurary" typo -> "January"
$ws.Range("B12").Value = "January"

# Add instructional text in the previously empty row 7 (no row insertion,
# rows 8-14 keep their positions)
$ws.Range("B7").Value = "select B8:B14 and drag"

# Update the active selection to B8
$ws.Range("B8").Select()
